$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) for rows 2-11
$ws.Range("C2").Value = -1.254546044999286
$ws.Range("D2").Value = 0.2181996216602013

$ws.Range("C3").Value = -0.07852526377260828
$ws.Range("D3").Value = 0.9378705727543892

$ws.Range("C4").Value = 2.279583853439159
$ws.Range("D4").Value = 0.02902757964453562

$ws.Range("C5").Value = 0.7871570161887943
$ws.Range("D5").Value = 0.4366429547782857

$ws.Range("C6").Value = 0.9935572529592049
$ws.Range("D6").Value = 0.3274572274493013

$ws.Range("C7").Value = 4.443285507125025
$ws.Range("D7").Value = 0.00008939871697943857

$ws.Range("C8").Value = 1.914359968027828
$ws.Range("D8").Value = 0.06402064744300362

$ws.Range("C9").Value = 2.518119980148215
$ws.Range("D9").Value = 0.01667219667356412

$ws.Range("C10").Value = 1.30584344453947
$ws.Range("D10").Value = 0.2003724854532902

$ws.Range("C11").Value = -1.890298525990514
$ws.Range("D11").Value = 0.06726686022743644
